$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1 previously held the "STATUS" header with a style; it now holds a
# trace/status message instead and loses its style formatting.
$ws.Range("C1").ClearFormats()
$ws.Range("C1").Value = "`tLogin Fail: Sorry this doesn't look like a valid email"

# Add per-row trace/status messages in column C for each account row.
$ws.Range("C2").Value = "`tLogin Successful"
$ws.Range("C3").Value = "`tLogin Fail: The password you entered is incorrect."
$ws.Range("C4").Value = "`tLogin Fail: The password you entered is incorrect."
$ws.Range("C5").Value = "`tLogin Fail: Hi, Cristian Camilo Isaza"
$ws.Range("C6").Value = "`tLogin Fail: Hi, Cristian Camilo Isaza"
